$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text (General/Text) representation
# rather than being auto-converted to numbers by Excel when the new
# value looks numeric (e.g. "1.00", "0.0611").

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.838.88'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.47%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.601.66'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.03%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '209.09'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.00'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.29%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.481'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -4.50%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.91%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0611'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.99%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '17.94'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0784'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.822.55'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.595.67'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.73%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.510'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.84%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.806.09'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.65%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.59'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.89%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0718'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.31%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.00'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '189.81'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.58%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.22%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -3.69%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.55%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.13%  '
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = 'Stellar'
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.128'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.38%  '
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '141.00'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.82%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.50%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.52'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -4.68%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.98%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.08%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.40%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -4.42%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.50%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.096.47'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.36%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.74%  '
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.794'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -8.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0151'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.07%  '
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'ImmutableX'
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.499'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.15%  '
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'Quant'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '95.71'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.69%  '
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'RocketPoolETH'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.735.39'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.10%  '
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.07'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.21%  '
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.744'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -4.33%  '
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0₆0112'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.10%  '
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '53.28'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.46%  '
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0512'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.02%  '
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.43'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.65%  '
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.409'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.24%  '
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'USDD'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.36%  '
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.28'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.27%  '
